$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '59.232.22'
$c.Style = $s
$ws.Range("E2").Value = '  -2.30%  '

$c = $ws.Range("D3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.576.24'
$c.Style = $s
$ws.Range("E3").Value = '  -2.53%  '

$ws.Range("E4").Value = '  +0.15%  '

$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '555.87'
$c.Style = $s
$ws.Range("E5").Value = '  -2.23%  '

$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '141.70'
$c.Style = $s
$ws.Range("E6").Value = '  -3.11%  '

$ws.Range("E7").Value = '  +0.07%  '

$c = $ws.Range("D8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.598'
$c.Style = $s
$ws.Range("E8").Value = '  -1.79%  '

$c = $ws.Range("D9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.581.28'
$c.Style = $s
$ws.Range("E9").Value = '  -2.82%  '

$ws.Range("E10").Value = '  -2.68%  '

$ws.Range("E11").Value = '  -1.03%  '

$ws.Range("E12").Value = '  +12.32%  '

$ws.Range("E13").Value = '  +2.41%  '

$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.027.26'
$c.Style = $s
$ws.Range("E14").Value = '  -3.00%  '

$c = $ws.Range("D15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '59.188.21'
$c.Style = $s
$ws.Range("E15").Value = '  -2.31%  '

$c = $ws.Range("D16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '23.02'
$c.Style = $s
$ws.Range("E16").Value = '  +4.22%  '

$c = $ws.Range("D17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0000137'
$c.Style = $s
$ws.Range("E17").Value = '  -1.43%  '

$c = $ws.Range("D18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.579.81'
$c.Style = $s
$ws.Range("E18").Value = '  -2.62%  '

$ws.Range("E19").Value = '  +0.11%  '

$c = $ws.Range("D20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '337.04'
$c.Style = $s
$ws.Range("E20").Value = '  -1.75%  '

$c = $ws.Range("D21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '10.35'
$c.Style = $s
$ws.Range("E21").Value = '  -1.05%  '

$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.41'
$c.Style = $s
$ws.Range("E22").Value = '  +0.34%  '

$c = $ws.Range("D23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = $s
$ws.Range("E23").Value = '  -0.19%  '

$c = $ws.Range("D24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '62.89'
$c.Style = $s
$ws.Range("E24").Value = '  -5.15%  '

$c = $ws.Range("D25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.471'
$c.Style = $s
$ws.Range("E25").Value = '  +7.26%  '

$ws.Range("E26").Value = '  +0.55%  '

$ws.Range("E27").Value = '  -3.36%  '

$c = $ws.Range("D28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.43'
$c.Style = $s
$ws.Range("E28").Value = '  +0.07%  '

$ws.Range("E29").Value = '  -5.07%  '

$ws.Range("E30").Value = '  +0.02%  '

$c = $ws.Range("D31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.19'
$c.Style = $s
$ws.Range("E31").Value = '  +0.19%  '

$c = $ws.Range("D32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.67'
$c.Style = $s
$ws.Range("E32").Value = '  -2.41%  '

$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '158.41'
$c.Style = $s
$ws.Range("E33").Value = '  -0.90%  '

$c = $ws.Range("D34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '19.08'
$c.Style = $s
$ws.Range("E34").Value = '  -0.64%  '

$ws.Range("E35").Value = '  -1.88%  '

$ws.Range("E36").Value = '  +0.48%  '

$c = $ws.Range("D37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.900'
$c.Style = $s
$ws.Range("E37").Value = '  +1.08%  '

$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '37.23'
$c.Style = $s
$ws.Range("E38").Value = '  -0.59%  '

$c = $ws.Range("D39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.853'
$c.Style = $s

$c = $ws.Range("D40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.48'
$c.Style = $s
$ws.Range("E40").Value = '  -2.83%  '

$ws.Range("E41").Value = '  +0.27%  '

$c = $ws.Range("D42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '291.08'
$c.Style = $s
$ws.Range("E42").Value = '  -3.01%  '

$c = $ws.Range("D43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '135.50'
$c.Style = $s
$ws.Range("E43").Value = '  +5.16%  '

$ws.Range("E44").Value = '  +0.17%  '

$c = $ws.Range("D45")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0976'
$c.Style = $s
$ws.Range("E45").Value = '  -1.09%  '

$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.591'
$c.Style = $s
$ws.Range("E46").Value = '  -2.00%  '

$ws.Range("E47").Value = '  -0.46%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0234'
$c.Style = $s
$ws.Range("E48").Value = '  -0.65%  '

$ws.Range("B49").Value = 'Hedera'
$ws.Range("C49").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0529'
$c.Style = $s
$ws.Range("E49").Value = '  -3.06%  '

$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '18.65'
$c.Style = $s
$ws.Range("E50").Value = '  -0.70%  '

$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.937.03'
$c.Style = $s
$ws.Range("E51").Value = '  -1.24%  '
